$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 covers the "Center on Brain Injury Research and Training" webinar.
# The conference-name cell previously also embedded "Webinar" in its text,
# and the location cell held a lowercase "webinar " placeholder. Split the
# title cleanly and move the properly-cased "Webinar " label into the
# conference-location column so the CV is ready to publish.
$ws.Range("N4").Value = "Center on Brain Injury Research and Training"
$ws.Range("P4").Value = "Webinar "

# Match the active selection that was in place when the workbook was last
# saved.
$ws.Range("P5").Select()
